# Applies the gh-pages data-refresh edit described by the commit
# "Update gh-pages to output generated at 456a3b4".
#
# Sheet layout (workbook.xml order):
#   1 = 展览     (Exhibition)
#   2 = 演出     (Performance)
#   3 = 本地生活 (Local life)      -- untouched
#   4 = 全部类型 (All categories)  -- combined listing, same rows reordered

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1 - 展览 (Exhibition): "want to go" counter bumps only
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value  = 400
$ws1.Range("F4").Value  = 1360
$ws1.Range("F5").Value  = 253
$ws1.Range("F6").Value  = 2622
$ws1.Range("F7").Value  = 987
$ws1.Range("F8").Value  = 19082
$ws1.Range("F10").Value = 2091
$ws1.Range("F11").Value = 703
$ws1.Range("F13").Value = 376
$ws1.Range("F14").Value = 640
$ws1.Range("F16").Value = 226
$ws1.Range("F20").Value = 227

# ---------------------------------------------------------------------
# Sheet 2 - 演出 (Performance): "want to go" counter bump only
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F16").Value = 79

# ---------------------------------------------------------------------
# Sheet 3 - 本地生活 (Local life): no changes
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# Sheet 4 - 全部类型 (All categories)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# "want to go" counter bumps on rows untouched by the row-25..30 shuffle
$ws4.Range("F6").Value  = 400
$ws4.Range("F8").Value  = 1360
$ws4.Range("F10").Value = 253
$ws4.Range("F13").Value = 2622
$ws4.Range("F14").Value = 987
$ws4.Range("F15").Value = 19082
$ws4.Range("F22").Value = 2091
$ws4.Range("F23").Value = 703
$ws4.Range("F36").Value = 227
$ws4.Range("F37").Value = 79

# The cancelled "凹凸世界ONLY" listing (old row 25) was dropped from this
# consolidated view, so every subsequent event in the 10-06..10-19 block
# shifts up one row, and a freshly-cancelled listing is appended at the
# row the block used to end on (row 30).

# Row 25: was 凹凸世界ONLY（取消） -> becomes 樱漫潮玩动漫游戏嘉年华
$ws4.Range("C25").Value = "广州·樱漫潮玩动漫游戏嘉年华"
$ws4.Range("D25").Value = "雄峰城B4座 广州番禺展览中心"
$ws4.Range("F25").Value = 376
$ws4.Range("G25").Value = 39.9
$ws4.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=91453"
$ws4.Range("I25").Value = "//i1.hdslb.com/bfs/openplatform/202409/xWacmeGY1726643680058.jpeg"

# Row 26: was 樱漫潮玩动漫游戏嘉年华 -> becomes 第五人格同人only4.0
$ws4.Range("C26").Value = "广州·第五人格同人only4.0"
$ws4.Range("D26").Value = "会江路巨大产业园5栋2楼 国际会议中心"
$ws4.Range("F26").Value = 640
$ws4.Range("G26").Value = 54
$ws4.Range("H26").Value = "https://show.bilibili.com/platform/detail.html?id=92168"
$ws4.Range("I26").Value = "//i0.hdslb.com/bfs/openplatform/202409/twg9GYOJ1725616149211.jpeg"

# Row 27: was 第五人格同人only4.0 -> becomes 第十届萌物语动漫嘉年华
# (force text format first so the date-shaped string isn't coerced to a
#  date serial number -- the source file stores these as plain text)
$ws4.Range("B27").NumberFormat = "@"
$ws4.Range("B27").Value = "2024-10-07"
$ws4.Range("C27").Value = "广州·第十届萌物语动漫嘉年华"
$ws4.Range("D27").Value = "雄峰城B4座 广州番禺展览中心"
$ws4.Range("E27").Value = "2024.10.07 10:00-10.07 17:00"
$ws4.Range("F27").Value = 213
$ws4.Range("G27").Value = 39.9
$ws4.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=91162"
$ws4.Range("I27").Value = "//i0.hdslb.com/bfs/openplatform/202409/cko2MQ2g1726643340997.jpeg"

# Row 28: was 第十届萌物语动漫嘉年华 -> becomes BanG Dream 同人Only
$ws4.Range("B28").NumberFormat = "@"
$ws4.Range("B28").Value = "2024-10-13"
$ws4.Range("C28").Value = "广州·BanG Dream 同人Only"
$ws4.Range("D28").Value = "华观路1932号 智慧城广百广场"
$ws4.Range("E28").Value = "2024.10.13 10:00-10.13 18:00"
$ws4.Range("F28").Value = 226
$ws4.Range("G28").Value = 58
$ws4.Range("H28").Value = "https://show.bilibili.com/platform/detail.html?id=92314"
$ws4.Range("I28").Value = "//i0.hdslb.com/bfs/openplatform/202409/HHkN2uUe1726045216331.png"

# Row 29: was BanG Dream 同人Only -> becomes Luca Stricagnoli 指弹吉他音乐会
$ws4.Range("B29").NumberFormat = "@"
$ws4.Range("B29").Value = "2024-10-19"
$ws4.Range("C29").Value = "广州·Luca Stricagnoli 2024《进化时间》指弹吉他音乐会"
$ws4.Range("D29").Value = "恩宁路265号3层 MaoLivehouse(永庆坊店)"
$ws4.Range("E29").Value = "2024.10.19 19:30-10.19 21:00"
$ws4.Range("F29").Value = 9
$ws4.Range("G29").Value = 220
$ws4.Range("H29").Value = "https://show.bilibili.com/platform/detail.html?id=91352"
$ws4.Range("I29").Value = "//i1.hdslb.com/bfs/openplatform/202408/ArhAlkP41724743278046.jpeg"

# Row 30: was Luca Stricagnoli 指弹吉他音乐会 -> becomes 次元喵喵动漫嘉年华02（取消）
$ws4.Range("C30").Value = "广州·次元喵喵动漫嘉年华02（取消）"
$ws4.Range("D30").Value = "东沙大道16号 广州健康方舟"
$ws4.Range("E30").Value = "2024.10.19 10:00-10.19 18:00"
$ws4.Range("F30").Value = 75
$ws4.Range("G30").Value = "不可售"
$ws4.Range("H30").Value = "https://show.bilibili.com/platform/detail.html?id=91566"
$ws4.Range("I30").Value = "//i1.hdslb.com/bfs/openplatform/202408/VJ9w4T6W1724046324480.jpeg"
